$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    # Force the cell to keep the literal text representation instead of
    # Excel auto-converting numeric-looking strings into numbers.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "24.375.40"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "1.648.85"
$ws.Range("E3").Value = "  -3.58%  "
Set-TextValue $ws.Range("D4") "1.005"
$ws.Range("E4").Value = "  +0.18%  "
Set-TextValue $ws.Range("D5") "311.05"
$ws.Range("E5").Value = "  -0.25%  "
Set-TextValue $ws.Range("D6") "1.002"
$ws.Range("E6").Value = "  +0.33%  "
Set-TextValue $ws.Range("D7") "0.3637"
$ws.Range("E7").Value = "  -3.41%  "
Set-TextValue $ws.Range("D8") "46.76"
$ws.Range("E8").Value = "  -5.83%  "
Set-TextValue $ws.Range("D9") "0.3230"
$ws.Range("E9").Value = "  -6.47%  "
Set-TextValue $ws.Range("D10") "1.115"
$ws.Range("E10").Value = "  -8.11%  "
Set-TextValue $ws.Range("D11") "0.06986"
$ws.Range("E11").Value = "  -7.61%  "
Set-TextValue $ws.Range("D12") "1.003"
$ws.Range("E12").Value = "  +0.37%  "
Set-TextValue $ws.Range("D13") "5.911"
$ws.Range("E13").Value = "  -6.67%  "
Set-TextValue $ws.Range("D14") "19.25"
$ws.Range("E14").Value = "  -9.56%  "
Set-TextValue $ws.Range("D15") "6.548"
$ws.Range("E15").Value = "  -7.59%  "
$ws.Range("D16").Value = "1.646.97"
$ws.Range("E16").Value = "  -3.76%  "
Set-TextValue $ws.Range("D17") "0.00001031"
$ws.Range("E17").Value = "  -9.30%  "
Set-TextValue $ws.Range("D18") "0.06586"
$ws.Range("E18").Value = "  -2.20%  "
Set-TextValue $ws.Range("D19") "1.001"
$ws.Range("E19").Value = "  +0.21%  "
Set-TextValue $ws.Range("D20") "77.74"
$ws.Range("E20").Value = "  -8.72%  "
Set-TextValue $ws.Range("D21") "5.894"
$ws.Range("E21").Value = "  -8.23%  "
Set-TextValue $ws.Range("D22") "15.48"
$ws.Range("E22").Value = "  -11.00%  "
Set-TextValue $ws.Range("D23") "12.42"
$ws.Range("E23").Value = "  -6.71%  "
$ws.Range("D24").Value = "24.355.82"
$ws.Range("E24").Value = "  -2.03%  "
Set-TextValue $ws.Range("D25") "2.480"
$ws.Range("E25").Value = "  +1.08%  "
Set-TextValue $ws.Range("D26") "2.287"
$ws.Range("E26").Value = "  -18.66%  "
Set-TextValue $ws.Range("D27") "145.30"
$ws.Range("E27").Value = "  -4.47%  "
Set-TextValue $ws.Range("D28") "18.44"
$ws.Range("E28").Value = "  -10.41%  "
$ws.Range("D29").Value = "1.830.80"
$ws.Range("E29").Value = "  -3.69%  "
Set-TextValue $ws.Range("D30") "123.15"
$ws.Range("E30").Value = "  -7.62%  "
Set-TextValue $ws.Range("D31") "1.168"
$ws.Range("E31").Value = "  -6.78%  "
Set-TextValue $ws.Range("D32") "4.054"
$ws.Range("E32").Value = "  -4.88%  "
Set-TextValue $ws.Range("D33") "5.594"
$ws.Range("E33").Value = "  -19.61%  "
Set-TextValue $ws.Range("D34") "0.08414"
$ws.Range("E34").Value = "  -4.85%  "
Set-TextValue $ws.Range("D35") "1.668"
$ws.Range("E35").Value = "  -7.04%  "
Set-TextValue $ws.Range("D36") "12.06"
$ws.Range("E36").Value = "  -13.76%  "
Set-TextValue $ws.Range("D37") "5.130"
$ws.Range("E37").Value = "  -8.97%  "
Set-TextValue $ws.Range("D38") "1.237"
$ws.Range("E38").Value = "  -3.85%  "
Set-TextValue $ws.Range("D39") "0.05956"
$ws.Range("E39").Value = "  -11.21%  "
Set-TextValue $ws.Range("D40") "0.02204"
$ws.Range("E40").Value = "  -9.07%  "
Set-TextValue $ws.Range("D41") "0.2039"
$ws.Range("E41").Value = "  -9.11%  "
Set-TextValue $ws.Range("D42") "8.047"
$ws.Range("E42").Value = "  -14.55%  "
$ws.Range("E43").Value = "  +0.24%  "
Set-TextValue $ws.Range("D44") "0.5837"
$ws.Range("E44").Value = "  -9.83%  "
Set-TextValue $ws.Range("D45") "3.753"
$ws.Range("E45").Value = "  -2.01%  "
Set-TextValue $ws.Range("D46") "12.51"
$ws.Range("E46").Value = "  -10.85%  "
Set-TextValue $ws.Range("D47") "0.5550"
$ws.Range("E47").Value = "  -10.28%  "
Set-TextValue $ws.Range("D48") "121.52"
$ws.Range("E48").Value = "  -6.93%  "
Set-TextValue $ws.Range("D49") "1.932"
$ws.Range("E49").Value = "  -9.73%  "
Set-TextValue $ws.Range("D50") "0.06877"
$ws.Range("E50").Value = "  -6.07%  "
Set-TextValue $ws.Range("D51") "1.173"
$ws.Range("E51").Value = "  -4.44%  "
